# refatorando o consolidador para modelo ETL
# Updates the absenteeism data rows 2-11 with the new ETL-sourced values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 64865
$ws.Range("B2").Value = "Gabriel Dias"
$ws.Range("C2").Value = "Vendas"
$ws.Range("D2").Value = "Viagem de negócios"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45082
$ws.Range("G2").Value = 7322.92

# Row 3
$ws.Range("A3").Value = 42332
$ws.Range("B3").Value = "Catarina Gonçalves"
$ws.Range("C3").Value = "P&D"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45092
$ws.Range("G3").Value = 5637.8

# Row 4 (Departamento/Motivo/Horas unchanged)
$ws.Range("A4").Value = 95632
$ws.Range("B4").Value = "Daniel Monteiro"
$ws.Range("F4").Value = 45085
$ws.Range("G4").Value = 10315.71

# Row 5 (Departamento unchanged)
$ws.Range("A5").Value = 70871
$ws.Range("B5").Value = "Luiz Felipe Alves"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 45081
$ws.Range("G5").Value = 10324.24

# Row 6
$ws.Range("A6").Value = 65710
$ws.Range("B6").Value = "Rafael Viana"
$ws.Range("C6").Value = "Engenharia"
$ws.Range("D6").Value = "Consulta médica"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 45081
$ws.Range("G6").Value = 11364.33

# Row 7
$ws.Range("A7").Value = 90300
$ws.Range("B7").Value = "João Lucas Costa"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 45081
$ws.Range("G7").Value = 6893.49

# Row 8
$ws.Range("A8").Value = 11564
$ws.Range("B8").Value = "João Vitor Dias"
$ws.Range("C8").Value = "Engenharia"
$ws.Range("D8").Value = "Doença"
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 45080
$ws.Range("G8").Value = 4322.12

# Row 9
$ws.Range("A9").Value = 29652
$ws.Range("B9").Value = "Pietra Nunes"
$ws.Range("C9").Value = "Atendimento ao Cliente"
$ws.Range("D9").Value = "Viagem de negócios"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45085
$ws.Range("G9").Value = 6611.15

# Row 10
$ws.Range("A10").Value = 81802
$ws.Range("B10").Value = "Vicente Pires"
$ws.Range("C10").Value = "TI"
$ws.Range("D10").Value = "Consulta médica"
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 45090
$ws.Range("G10").Value = 10785.18

# Row 11
$ws.Range("A11").Value = 61251
$ws.Range("B11").Value = "Kevin Aragão"
$ws.Range("C11").Value = "Vendas"
$ws.Range("D11").Value = "Viagem de negócios"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45099
$ws.Range("G11").Value = 8604.17
